# Rows 69-82 on the "Artfynd" sheet got re-sorted: the species-related
# columns (A, B, D, E, F, G, H, Q, R) for each of those rows now hold the
# values that used to belong to a different row within the same block,
# while every other column (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AT, AW, AX, AY, ...) stays put. Below is the row->row mapping of
# "new row N gets the data that used to live in row M" (both are within
# 69..82, and the mapping is a permutation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    69 = 74
    70 = 82
    71 = 80
    72 = 77
    73 = 72
    74 = 78
    75 = 79
    76 = 75
    77 = 73
    78 = 76
    79 = 81
    80 = 71
    81 = 70
    82 = 69
}

# Columns whose values move together with the row mapping.
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18)   # A, B, D, E, F, G, H, Q, R

# Snapshot all the source values first so that the column values of one
# row are not clobbered before being used as the source for another row.
# NOTE: use .Value2 for reads (this runtime's .Value getter does not
# reliably resolve through property reflection), .Value for writes.
$snapshot = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowValues = @{}
        foreach ($col in $cols) {
            $rowValues[$col] = $ws.Cells.Item($srcRow, $col).Value2
        }
        $snapshot[$srcRow] = $rowValues
    }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value = $rowValues[$col]
    }
}
